# Updated 7 test cases for DC Unit scripts with new loading method details
#
# Adds a new "DC Unit Loading Details Name" column (with "Current (DC Units)"
# and "Current (worst case)" rows) to the "Add Devices Loop A", "Add Devices
# Loop B" and "Panel LED" sheets, and touches the page setup on the three
# "no new data" sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: Add Devices Loop A  (new column C)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Add Devices Loop A")

# Values first, in the same order the original author typed them, so the
# shared-string table ends up with the same index assignment (51/52/53).
$ws1.Range("C2").Value = "Current (DC Units)"
$ws1.Range("C1").Value = "DC Unit Loading Details Name"
$ws1.Range("C3").Value = "Current (worst case)"

# Formats: C1 matches the bold "table header" style used by row 7 (A7:N7);
# C2/C3 match the "table value" style used by row 8 (A8:N8).
$ws1.Range("A7").Copy() | Out-Null
$ws1.Range("C1").PasteSpecial(-4122) | Out-Null
$ws1.Range("A8").Copy() | Out-Null
$ws1.Range("C2:C3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws1.Columns.Item(3).ColumnWidth = 26.33203125

$ws1.Range("C2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet: Add Devices Loop B  (new column C, same text/style, narrower col)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Add Devices Loop B")

$ws2.Range("C1").Value = "DC Unit Loading Details Name"
$ws2.Range("C2").Value = "Current (DC Units)"
$ws2.Range("C3").Value = "Current (worst case)"

$ws2.Range("A7").Copy() | Out-Null
$ws2.Range("C1").PasteSpecial(-4122) | Out-Null
$ws2.Range("A8").Copy() | Out-Null
$ws2.Range("C2:C3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws2.Range("C1:C3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet: Panel LED  (new column F, plus a matching width on column C)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Panel LED")

$ws3.Range("F1").Value = "DC Unit Loading Details Name"
$ws3.Range("F2").Value = "Current (DC Units)"
$ws3.Range("F3").Value = "Current (worst case)"

# Panel LED's own rows don't carry the s=8/s=12 header/value style pair, so
# pull the formats from "Add Devices Loop A", which does.
$ws1.Range("A7").Copy() | Out-Null
$ws3.Range("F1").PasteSpecial(-4122) | Out-Null
$ws1.Range("A8").Copy() | Out-Null
$ws3.Range("F2:F3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws3.Columns.Item(3).ColumnWidth = 26.33203125
$ws3.Columns.Item(6).ColumnWidth = 26.33203125

$ws3.PageSetup.Orientation = 1

$ws3.Range("F1:F3").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheets: Delete Devices Loop A / B -- page setup touched, no data changes
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Delete Devices Loop A")
$ws4.PageSetup.Orientation = 1

$ws5 = $wb.Worksheets.Item("Delete Devices Loop B")
$ws5.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# Restore the originally active sheet/selection
# ---------------------------------------------------------------------------
$ws1.Activate()
